$wb = $excel.ActiveWorkbook

$wsOv = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Overview sheet: append row 4 for the new handback file
# ---------------------------------------------------------------------------
$loOv = $wsOv.ListObjects.Item(1)
$loOv.ListRows.Add()

$wsOv.Range("A4").Value = "f21c113d-d635-4905-a529-ed7823725a01.md"
$wsOv.Range("C4").Value = ".md"
$wsOv.Range("E4").Value = "Handed back: in sync with en-US"
$wsOv.Range("F4").Value = "Handed back: in sync with en-US"
$wsOv.Range("G4").Value = "2016-09-02 12:51:25"
$wsOv.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOv.Hyperlinks.Add($wsOv.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/HEAD/e2e/f21c113d-d635-4905-a529-ed7823725a01.md", "", "", "e2e\f21c113d-d635-4905-a529-ed7823725a01.md")

# ---------------------------------------------------------------------------
# zh-cn sheet: append row 4
# ---------------------------------------------------------------------------
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add()

$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Handed back: in sync with en-US"
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "'True"
$wsZh.Range("G4").Value = "f21c113d-d635-4905-a529-ed7823725a01.03d66b3d71282592a133485428355c1063cb0e9d.zh-cn.xlf"
$wsZh.Range("H4").Value = "2016-09-02 12:51:20"
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("J4").Value = "f21c113d-d635-4905-a529-ed7823725a01.03d66b3d71282592a133485428355c1063cb0e9d.zh-cn.xlf"
$wsZh.Range("K4").Value = "2016-09-02 12:51:37"
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L4").Value = "'"
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("N4").Value = "'"
$wsZh.Range("O4").Value = "'False"
$wsZh.Range("P4").Value = "'"

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/HEAD/e2e/f21c113d-d635-4905-a529-ed7823725a01.md", "", "", "f21c113d-d635-4905-a529-ed7823725a01.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/HEAD/e2e/f21c113d-d635-4905-a529-ed7823725a01.md", "", "", "f21c113d-d635-4905-a529-ed7823725a01.md")

# ---------------------------------------------------------------------------
# de-de sheet: append row 4
# ---------------------------------------------------------------------------
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add()

$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "'True"
$wsDe.Range("G4").Value = "f21c113d-d635-4905-a529-ed7823725a01.03d66b3d71282592a133485428355c1063cb0e9d.de-de.xlf"
$wsDe.Range("H4").Value = "2016-09-02 12:51:25"
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("J4").Value = "f21c113d-d635-4905-a529-ed7823725a01.03d66b3d71282592a133485428355c1063cb0e9d.de-de.xlf"
$wsDe.Range("K4").Value = "2016-09-02 12:51:45"
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L4").Value = "'"
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("N4").Value = "'"
$wsDe.Range("O4").Value = "'False"
$wsDe.Range("P4").Value = "'"

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/HEAD/e2e/f21c113d-d635-4905-a529-ed7823725a01.md", "", "", "f21c113d-d635-4905-a529-ed7823725a01.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/HEAD/e2e/f21c113d-d635-4905-a529-ed7823725a01.md", "", "", "f21c113d-d635-4905-a529-ed7823725a01.md")

Write-Host "Handback report row added for f21c113d-d635-4905-a529-ed7823725a01"
